$d = $word.ActiveDocument

# 1. Right-align the "姓名：... 日期：..." paragraph (2nd body paragraph).
$d.Paragraphs.Item(2).Alignment = 2

# 2. Apply the "Table Grid" style to the sign-in table and center it on
#    the page.
$t = $d.Tables.Item(1)
$t.Style = "Table Grid"
$t.Alignment = 1

# 3. Center every paragraph inside every cell of the table.
foreach ($row in $t.Rows) {
    foreach ($cell in $row.Cells) {
        foreach ($p in $cell.Range.Paragraphs) {
            $p.Alignment = 1
        }
    }
}
